# Updates the "Estado de Cuenta" debtor table: previous account-statement
# rows are replaced by the new data, regrouped per worker (most recent
# period first) as produced by the refreshed database macro.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$people = @(
    @{ Doc = "73122753";   Name = "DAVID MORELO BENITEZ" },
    @{ Doc = "9146245";    Name = "VICTOR MANUEL MORELO BENITEZ" },
    @{ Doc = "73202100";   Name = "JUAN CARLOS VALLEJO CASTRO" },
    @{ Doc = "1004346188"; Name = "PEREZ RODRIGUEZ ALEPH DALEPH" }
)

$periods = @(
    @{ Period = "2104"; Value = 30430 },
    @{ Period = "2103"; Value = 35112 },
    @{ Period = "2102"; Value = 35112 },
    @{ Period = "2101"; Value = 35112 },
    @{ Period = "2012"; Value = 35112 },
    @{ Period = "2011"; Value = 35112 },
    @{ Period = "2010"; Value = 35112 },
    @{ Period = "2009"; Value = 35112 },
    @{ Period = "2008"; Value = 15215 }
)

$row = 16
foreach ($person in $people) {
    foreach ($p in $periods) {
        $ws.Range("B$row").Value = "CC"
        $ws.Range("C$row").Value = $person.Doc
        $ws.Range("D$row").Value = $person.Name
        $ws.Range("E$row").Value = $p.Period
        $ws.Range("F$row").Value = $p.Value
        $ws.Range("G$row").Value = 877803
        $row++
    }
}
